$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily refresh of the cryptos price/volume table (GitHub Actions job).
# Price cells in column D are plain text (values like "3.137.48" use
# "." as a thousands separator and are not valid Excel numbers), so we
# force a Text format before writing to stop Excel from silently
# reinterpreting numeric-looking strings (e.g. "5.90") as numbers, then
# restore the original "Normal" style so formatting is left untouched.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.561.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.128.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.93%  '
$ws.Range("E7").Value = '  +6.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.374'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.42%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.127.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.763'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.95%  '
$ws.Range("E12").Value = '  +3.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.213.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.706.70'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.136.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  +3.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000215'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '452.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.47%  '
$ws.Range("E25").Value = '  -1.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '92.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("E27").Value = '  -3.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.293.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.181'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +12.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.239'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +18.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.122'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +40.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.81%  '
$ws.Range("E34").Value = '  +36.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.165'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.22'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +26.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '502.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.57%  '
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.424'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.701'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.28%  '
